# Update "Max Measurement P-Value" (column G) values for several rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.043
$ws.Range("G3").Value = 0.002
$ws.Range("G5").Value = 0.002
$ws.Range("G6").Value = 0.001

# Remove the now-unneeded "Max Change Rate (Normalized)" column (H) entirely,
# which also drops its now-unused shared string and shrinks the sheet dimension.
$ws.Range("H1").EntireColumn.Delete()
